# Commit: Sat, May 16, 2020  6:04:50 AM
#
# 1) Slide 16's table switches from the custom "Table_0" style
#    ({79D18943-65EC-4B84-A353-3A86E70599C5}, defined in ppt/tableStyles.xml)
#    to the built-in PowerPoint table style {73DD70DF-311F-4851-9BDF-A13756CC3B33}.
#
# 2) The deck's theme (ppt/theme/theme2.xml, the one actually driving the
#    slide master / slides) is recoloured from the "Integral" palette to the
#    stock Office palette (the colours that used to live, unused, in
#    ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation

# --- 1) table style on slide 16 -------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{73DD70DF-311F-4851-9BDF-A13756CC3B33}")
    }
}

# --- 2) theme colour scheme -------------------------------------------------
# ThemeColorScheme.Colors(n).RGB uses the classic COLORREF (0x00BBGGRR)
# ordering, so each target "RRGGBB" hex value has to be byte-reversed.
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

$scheme.Colors(1).RGB  = 0x000000  # dk1      000000
$scheme.Colors(2).RGB  = 0xFFFFFF  # lt1      FFFFFF
$scheme.Colors(3).RGB  = 0x6A5444  # dk2      44546A
$scheme.Colors(4).RGB  = 0xE6E6E7  # lt2      E7E6E6
$scheme.Colors(5).RGB  = 0xD59B5B  # accent1  5B9BD5
$scheme.Colors(6).RGB  = 0x317DED  # accent2  ED7D31
$scheme.Colors(7).RGB  = 0xA5A5A5  # accent3  A5A5A5
$scheme.Colors(8).RGB  = 0x00C0FF  # accent4  FFC000
$scheme.Colors(9).RGB  = 0xC47244  # accent5  4472C4
$scheme.Colors(10).RGB = 0x47AD70  # accent6  70AD47
$scheme.Colors(11).RGB = 0xC16305  # hlink    0563C1
$scheme.Colors(12).RGB = 0x724F95  # folHlink 954F72
